$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap content of column B through AD between row pairs (column A id stays fixed per row)
$rowA = $ws.Range("B87:AD87").Value2
$rowB = $ws.Range("B88:AD88").Value2
$ws.Range("B87:AD87").Value2 = $rowB
$ws.Range("B88:AD88").Value2 = $rowA

$rowA = $ws.Range("B130:AD130").Value2
$rowB = $ws.Range("B131:AD131").Value2
$ws.Range("B130:AD130").Value2 = $rowB
$ws.Range("B131:AD131").Value2 = $rowA

$rowA = $ws.Range("B183:AD183").Value2
$rowB = $ws.Range("B184:AD184").Value2
$ws.Range("B183:AD183").Value2 = $rowB
$ws.Range("B184:AD184").Value2 = $rowA

$rowA = $ws.Range("B205:AD205").Value2
$rowB = $ws.Range("B206:AD206").Value2
$ws.Range("B205:AD205").Value2 = $rowB
$ws.Range("B206:AD206").Value2 = $rowA

$rowA = $ws.Range("B216:AD216").Value2
$rowB = $ws.Range("B217:AD217").Value2
$ws.Range("B216:AD216").Value2 = $rowB
$ws.Range("B217:AD217").Value2 = $rowA

$rowA = $ws.Range("B226:AD226").Value2
$rowB = $ws.Range("B227:AD227").Value2
$ws.Range("B226:AD226").Value2 = $rowB
$ws.Range("B227:AD227").Value2 = $rowA

$rowA = $ws.Range("B235:AD235").Value2
$rowB = $ws.Range("B236:AD236").Value2
$ws.Range("B235:AD235").Value2 = $rowB
$ws.Range("B236:AD236").Value2 = $rowA

$rowA = $ws.Range("B240:AD240").Value2
$rowB = $ws.Range("B241:AD241").Value2
$ws.Range("B240:AD240").Value2 = $rowB
$ws.Range("B241:AD241").Value2 = $rowA

$rowA = $ws.Range("B259:AD259").Value2
$rowB = $ws.Range("B261:AD261").Value2
$ws.Range("B259:AD259").Value2 = $rowB
$ws.Range("B261:AD261").Value2 = $rowA

# 3-way rotation for rows 145, 146, 147: new145=old146, new146=old147, new147=old145
$row145 = $ws.Range("B145:AD145").Value2
$row146 = $ws.Range("B146:AD146").Value2
$row147 = $ws.Range("B147:AD147").Value2
$ws.Range("B145:AD145").Value2 = $row146
$ws.Range("B146:AD146").Value2 = $row147
$ws.Range("B147:AD147").Value2 = $row145
